$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.272107839584351
$ws.Range("B1").Value = 2.275051593780518
$ws.Range("C1").Value = 4.645138740539551
$ws.Range("D1").Value = 3.045556783676147
$ws.Range("E1").Value = 1.352442741394043
